$wb = $excel.ActiveWorkbook

$newPath = "C:\Katalon_mastercopy520\DataCommons_Automation\OutputFiles\TC05_Canine_Filter_Breed-Beagle_Neo4jData.xlsx"

$wsMessage = $wb.Worksheets.Item("Message")
$wsMessage.Range("A10").Value = $newPath

$wsCypherOutputMessage = $wb.Worksheets.Item("CypherOutput_Message")
$wsCypherOutputMessage.Range("A10").Value = $newPath

$wsStatOutputMessage = $wb.Worksheets.Item("StatOutput_Message")
$wsStatOutputMessage.Range("A20").Value = $newPath
